$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = 1
$ws.Cells.Item(110, 3).Value = "2024-06-17 06:18:43"
$ws.Cells.Item(110, 4).Value = 200
$ws.Cells.Item(110, 5).Value = 9

$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = 2
$ws.Cells.Item(111, 3).Value = "2024-06-17 06:18:44"
$ws.Cells.Item(111, 4).Value = 200
$ws.Cells.Item(111, 5).Value = 0
